# Auto-generated edit script: refreshes the cryptos price/volume table
# (GitHub Actions data refresh) plus a ShibaInu/WstETH row swap.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.735.16'
$ws.Range("E2").Value = '  +0.22%  '
$ws.Range("D3").Value = '2.615.66'
$ws.Range("E3").Value = '  -2.32%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.19'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.17'
$ws.Range("E6").Value = '  -1.21%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  -0.72%  '
$ws.Range("D9").Value = '2.617.88'
$ws.Range("E9").Value = '  -2.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.168'
$ws.Range("E10").Value = '  -0.58%  '
$ws.Range("E11").Value = '  +0.74%  '
$ws.Range("E12").Value = '  +1.51%  '
$ws.Range("E13").Value = '  -2.15%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '3.097.35'
$ws.Range("E14").Value = '  -2.16%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000188'
$ws.Range("E15").Value = '  +2.43%  '
$ws.Range("D16").Value = '71.721.41'
$ws.Range("E16").Value = '  +0.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.57'
$ws.Range("E17").Value = '  -2.06%  '
$ws.Range("D18").Value = '2.632.86'
$ws.Range("E18").Value = '  -1.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.03'
$ws.Range("E19").Value = '  -0.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.83'
$ws.Range("E20").Value = '  -1.93%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '372.52'
$ws.Range("E21").Value = '  +1.21%  '
$ws.Range("E22").Value = '  -2.00%  '
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.23'
$ws.Range("E24").Value = '  -0.48%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("E26").Value = '  -2.82%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.31'
$ws.Range("E27").Value = '  -5.34%  '
$ws.Range("D28").Value = '2.763.29'
$ws.Range("E28").Value = '  -1.76%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("E30").Value = '  -1.52%  '
$ws.Range("E31").Value = '  -1.50%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '484.29'
$ws.Range("E32").Value = '  -4.09%  '
$ws.Range("E33").Value = '  +2.03%  '
$ws.Range("E34").Value = '  -0.94%  '
$ws.Range("E35").Value = '  +0.21%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '159.78'
$ws.Range("E36").Value = '  -1.78%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.115'
$ws.Range("E37").Value = '  +6.44%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.23'
$ws.Range("E38").Value = '  -0.75%  '
$ws.Range("E39").Value = '  -0.74%  '
$ws.Range("E40").Value = '  -1.56%  '
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("E42").Value = '  -4.83%  '
$ws.Range("E43").Value = '  -0.81%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.82'
$ws.Range("E44").Value = '  -3.33%  '
$ws.Range("E45").Value = '  -2.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '38.91'
$ws.Range("E46").Value = '  -0.67%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '149.16'
$ws.Range("E47").Value = '  -3.98%  '
$ws.Range("E48").Value = '  -1.86%  '
$ws.Range("E49").Value = '  -1.52%  '
$ws.Range("E50").Value = '  -4.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.603'
$ws.Range("E51").Value = '  +0.27%  '
